# Auto-generated edit script applying numeric updates to Belias_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1680350.6
$ws.Range("J11").Value = 20526
$ws.Range("L11").Value = 20526
$ws.Range("N11").Value = -20814
$ws.Range("H13").Value = 5108.6665
$ws.Range("J13").Value = 5108.6665
$ws.Range("L13").Value = 5108.6665
$ws.Range("N13").Value = -5396.6665
$ws.Range("H28").Value = 20578.5
$ws.Range("I28").Value = 15867.75
$ws.Range("J28").Value = 30000
$ws.Range("K28").Value = 15867.75
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = -15675.75
$ws.Range("N28").Value = -30384
$ws.Range("H32").Value = 6494.2314
$ws.Range("I32").Value = 3127.5898
$ws.Range("K32").Value = 3127.5898
$ws.Range("M32").Value = -2840.5898
$ws.Range("H69").Value = 112450.125
$ws.Range("J69").Value = 112450.125
$ws.Range("L69").Value = 112450.125
$ws.Range("N69").Value = -113948.125
$ws.Range("H72").Value = 112450.125
$ws.Range("J72").Value = 112450.125
$ws.Range("L72").Value = 337350.375
$ws.Range("N72").Value = -344838.375
$ws.Range("H93").Value = 29224
$ws.Range("J93").Value = 29224
$ws.Range("L93").Value = 29224
$ws.Range("N93").Value = -34216
$ws.Range("H99").Value = 20578.5
$ws.Range("I99").Value = 15867.75
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 15867.75
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -12872.75
$ws.Range("N99").Value = -35990
$ws.Range("H132").Value = 1804.6735
$ws.Range("I132").Value = 1158.6207
$ws.Range("J132").Value = 2741.45
$ws.Range("K132").Value = 3475.8621
$ws.Range("L132").Value = 8224.349999999999
$ws.Range("M132").Value = -945.8620999999998
$ws.Range("N132").Value = -13284.35

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30472
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2268.9167
$ws.Range("I31").Value = 1239.4615
$ws.Range("J31").Value = 2651.2856
$ws.Range("K31").Value = 1239.4615
$ws.Range("L31").Value = 2651.2856
$ws.Range("M31").Value = -944.4614999999999
$ws.Range("N31").Value = -3241.2856
$ws.Range("H34").Value = 2268.9167
$ws.Range("I34").Value = 1239.4615
$ws.Range("J34").Value = 2651.2856
$ws.Range("K34").Value = 1239.4615
$ws.Range("L34").Value = 2651.2856
$ws.Range("M34").Value = -1037.4615
$ws.Range("N34").Value = -3055.2856
$ws.Range("H97").Value = 23595.2
$ws.Range("J97").Value = 23524.25
$ws.Range("L97").Value = 23524.25
$ws.Range("N97").Value = -25506.25
$ws.Range("H99").Value = 2063.8
$ws.Range("I99").Value = 1915
$ws.Range("J99").Value = 2361.4
$ws.Range("K99").Value = 1915
$ws.Range("L99").Value = 2361.4
$ws.Range("M99").Value = -417
$ws.Range("N99").Value = -5357.4
$ws.Range("H105").Value = 765.5625
$ws.Range("I105").Value = 871.2857
$ws.Range("J105").Value = 683.3333
$ws.Range("K105").Value = 871.2857
$ws.Range("L105").Value = 683.3333
$ws.Range("M105").Value = 875.7143
$ws.Range("N105").Value = -4177.3333
$ws.Range("H126").Value = 2063.8
$ws.Range("I126").Value = 1915
$ws.Range("J126").Value = 2361.4
$ws.Range("K126").Value = 5745
$ws.Range("L126").Value = 7084.200000000001
$ws.Range("M126").Value = -3275
$ws.Range("N126").Value = -12024.2
$ws.Range("H141").Value = 69911.89999999999
$ws.Range("J141").Value = 69911.89999999999
$ws.Range("L141").Value = 69911.89999999999
$ws.Range("N141").Value = -80271.89999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1615.3572
$ws.Range("I34").Value = 551
$ws.Range("J34").Value = 2206.6667
$ws.Range("K34").Value = 1653
$ws.Range("L34").Value = 6620.000100000001
$ws.Range("M34").Value = -1569
$ws.Range("N34").Value = -6788.000100000001
$ws.Range("H92").Value = 1143.8572
$ws.Range("J92").Value = 1250.75
$ws.Range("L92").Value = 3752.25
$ws.Range("N92").Value = -6248.25
$ws.Range("H109").Value = 4821.9165
$ws.Range("J109").Value = 5903.8887
$ws.Range("L109").Value = 17711.6661
$ws.Range("N109").Value = -19791.6661
$ws.Range("H131").Value = 881.14
$ws.Range("I131").Value = 445.66666
$ws.Range("J131").Value = 940.5227
$ws.Range("K131").Value = 1336.99998
$ws.Range("L131").Value = 2821.5681
$ws.Range("M131").Value = 3703.00002
$ws.Range("N131").Value = -12901.5681

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 20638.715
$ws.Range("I99").Value = 4235.5
$ws.Range("J99").Value = 27200
$ws.Range("K99").Value = 4235.5
$ws.Range("L99").Value = 27200
$ws.Range("M99").Value = -1989.5
$ws.Range("N99").Value = -31692
$ws.Range("H122").Value = 1283.2174
$ws.Range("I122").Value = 1271.7142
$ws.Range("J122").Value = 1404
$ws.Range("K122").Value = 3815.1426
$ws.Range("L122").Value = 4212
$ws.Range("M122").Value = -1365.1426
$ws.Range("N122").Value = -9112

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3289.2
$ws.Range("I7").Value = 2748
$ws.Range("K7").Value = 2748
$ws.Range("M7").Value = -2636
$ws.Range("H22").Value = 903.3889
$ws.Range("I22").Value = 764.1429000000001
$ws.Range("J22").Value = 992
$ws.Range("K22").Value = 764.1429000000001
$ws.Range("L22").Value = 992
$ws.Range("M22").Value = -469.1429000000001
$ws.Range("N22").Value = -1582
$ws.Range("H27").Value = 903.3889
$ws.Range("I27").Value = 764.1429000000001
$ws.Range("J27").Value = 992
$ws.Range("K27").Value = 764.1429000000001
$ws.Range("L27").Value = 992
$ws.Range("M27").Value = -657.1429000000001
$ws.Range("N27").Value = -1206
$ws.Range("H40").Value = 1557.9565
$ws.Range("I40").Value = 1448.4117
$ws.Range("K40").Value = 1448.4117
$ws.Range("M40").Value = -1312.4117
$ws.Range("H122").Value = 3456.4092
$ws.Range("I122").Value = 4280.125
$ws.Range("K122").Value = 12840.375
$ws.Range("M122").Value = -10390.375
$ws.Range("H126").Value = 3289.2
$ws.Range("I126").Value = 2748
$ws.Range("K126").Value = 8244
$ws.Range("M126").Value = -5774
$ws.Range("H131").Value = 30000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -40080
$ws.Range("H132").Value = 3187.1702
$ws.Range("I132").Value = 3072.6765
$ws.Range("J132").Value = 3486.6155
$ws.Range("K132").Value = 9218.029500000001
$ws.Range("L132").Value = 10459.8465
$ws.Range("M132").Value = -6688.029500000001
$ws.Range("N132").Value = -15519.8465

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 101826.086
$ws.Range("I122").Value = 1885.5714
$ws.Range("J122").Value = 241742.8
$ws.Range("K122").Value = 5656.7142
$ws.Range("L122").Value = 725228.3999999999
$ws.Range("M122").Value = -3206.7142
$ws.Range("N122").Value = -730128.3999999999
